$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 73; this shifts rows 73..162 down to 74..163,
# carrying their existing values/styles with them (matches the observed diff:
# every row from 73 onward absorbs the prior row's D/J/K/M/P data, and the old
# row 162 becomes the new row 163).
$ws.Rows("73:73").Insert()

# Populate the newly inserted row 73 with the new weekly data point.
$ws.Range("A73").Value = 8
$ws.Range("B73").Value = "Terminal La Palmera de La Serena"
$ws.Range("C73").Value = "Coquimbo"
$ws.Range("D73").Value = 44483
$ws.Range("E73").Value = 4
$ws.Range("F73").Value = 100112012
$ws.Range("G73").Value = "Espinaca"
$ws.Range("H73").Value = "Sin especificar"
$ws.Range("I73").Value = "Primera"
$ws.Range("J73").Value = 3060
$ws.Range("K73").Value = 400
$ws.Range("L73").Value = 500
$ws.Range("M73").Value = 450
$ws.Range("N73").Value = "$/atado 300 a 500 gramos"
$ws.Range("O73").Value = "Provincia del Elquí"
$ws.Range("P73").Value = 900
$ws.Range("Q73").Value = 0.5
$ws.Range("R73").Value = "Hortaliza"

# Match the date-formatted style used by the other D-column cells.
$ws.Range("D73").NumberFormat = "YYYY-MM-DD HH:MM:SS"
